# Fixed inconsistencies in BOM
#  - "Quatity" header typo corrected to "Quantity"
#  - R2 (row 10) quantity corrected from 2 to 1 (only one designator, R2)
#  - R5,R1 (row 13) quantity corrected from 4 to 2 (two designators, R5 and R1)
#  - cosmetic: column widths re-fit, active selection moved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo: "Quatity" -> "Quantity" ---
$ws.Range("B2").Value = "Quantity"

# --- Fix quantity data inconsistencies ---
$ws.Range("B10").Value = 1
$ws.Range("B13").Value = 2

# --- Cosmetic: column widths (re-autofit after text changes) ---
$ws.Columns.Item(1).ColumnWidth = 4.666666666666667
$ws.Columns.Item(2).ColumnWidth = 8.666666666666666
$ws.Columns.Item(3).ColumnWidth = 11.5
$ws.Columns.Item(4).ColumnWidth = 29.666666666666668
$ws.Columns.Item(5).ColumnWidth = 29.5
$ws.Columns.Item(6).ColumnWidth = 27.666666666666668
$ws.Columns.Item(7).ColumnWidth = 15.833333333333334

# --- Cosmetic: move active selection ---
[void]$ws.Range("D23").Select()
